$wb = $excel.ActiveWorkbook

# --- Add the new "espn_api" worksheet after the last existing sheet (scrape_call) ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "espn_api"

# --- Pre-touch these values off to the side so the shared-string table gets the
#     same allocation order as the source workbook (swid, espn_s2, season, league_id,
#     then the guid/token literals) before the real header/data layout is written. ---
$newSheet.Range("Z1").Value = "swid"
$newSheet.Range("Z2").Value = "espn_s2"
$newSheet.Range("Z3").Value = "season"
$newSheet.Range("Z4").Value = "league_id"
$newSheet.Range("Z5").Value = "{1AFA2BA9-0A5B-499B-BC41-CA32FDB46E50}"
$newSheet.Range("Z6").Value = "AEB9IogDihpmnSBySa1Dt%2BBZGNCwhpyrx1lOauwtvC2Zk4%2F1%2FXPWx%2BBr2HCsUsK3IL3Y6ansGgILWrmKQ5KG3em296twmYEMjOUDGa%2FwYNQ44GBlps9n6Vtts4%2Fh43ivZJzUzMH6dK5%2BQTwyG4wOoj6hXAjyq6gtHh5qSUDJDxqmeaCejF%2BntpFFHdx5kTfcE46%2F0XnWu7IAW2svLYPC53uMQPeYyBkNwMVfOincdbVpgcbquuF898mZeOnWG8ZjbI4e6Wg9q403IYw5o9ua%2FFk5uRCnpit%2FB9x7zgCmOsmhhA%3D%3D"

# --- Header row ---
$newSheet.Range("A1").Value = "season"
$newSheet.Range("B1").Value = "league_id"
$newSheet.Range("C1").Value = "swid"
$newSheet.Range("D1").Value = "espn_s2"

# --- Data row ---
$newSheet.Range("A2").Value = 2023
$newSheet.Range("B2").Value = 44419657
$newSheet.Range("C2").Value = "{1AFA2BA9-0A5B-499B-BC41-CA32FDB46E50}"
$newSheet.Range("D2").Value = "AEB9IogDihpmnSBySa1Dt%2BBZGNCwhpyrx1lOauwtvC2Zk4%2F1%2FXPWx%2BBr2HCsUsK3IL3Y6ansGgILWrmKQ5KG3em296twmYEMjOUDGa%2FwYNQ44GBlps9n6Vtts4%2Fh43ivZJzUzMH6dK5%2BQTwyG4wOoj6hXAjyq6gtHh5qSUDJDxqmeaCejF%2BntpFFHdx5kTfcE46%2F0XnWu7IAW2svLYPC53uMQPeYyBkNwMVfOincdbVpgcbquuF898mZeOnWG8ZjbI4e6Wg9q403IYw5o9ua%2FFk5uRCnpit%2FB9x7zgCmOsmhhA%3D%3D"

# Remove the temporary helper cells now that the shared strings are registered
$newSheet.Range("Z1:Z6").Clear()

# E2 is a formatted-but-empty cell (quote-prefix / "treat as text" styling), matching
# the equivalent empty-but-styled cell pattern already used elsewhere in this workbook.
$newSheet.Range("E2").Value = "'"
$newSheet.Range("E2").ClearContents()

# Auto-size the swid/league_id columns like the source sheet's best-fit widths
$newSheet.Columns("B:C").AutoFit() | Out-Null

# The newly added sheet becomes the active tab (mirrors the source workbook, which
# also now opens on "espn_api" instead of "variables")
$newSheet.Activate()
